$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7

$ws.Cells.Item($row, 1).Value = 6
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "17.04.2023 12:44 (CET)"
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = '{"ProposedVersion":"https://gitlab.intra.infineon.com/semantic-web-projects/digital-reference/order_management/-/commit/50873a28b7e97ba8d65492c7a7938f0fc336ac02","UpdatedVersion":"https://gitlab.intra.infineon.com/semantic-web-projects/digital-reference/order_management/-/commit/90c783a6e4ff122dbc160ccf0a2745e78a024c73","Domain":"Supply Chain","LobeOwner":"member1","Result":"accept by lobe owner"}'
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "3eabfa48eb39069846a1b161e3d2f19b23a539121e9fdd3692294e375171601c"

$ws.Range("A1:D7").Style = "Normal"
